$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: "8020 Joining plate 2rcx4" part is replaced by "inside corner connector",
# with a new cost-per-unit-size (12.62 -> 7) and number-at-size (16 -> 20).
$ws.Range("A5").Value = "inside corner connector"
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 20

# Row 6 ("tnuts"): E6's formula no longer depends on the row-10 "panel
# retainer" quantity (E10) - that dependency is replaced by a literal 10 -
# update this before deleting row 10 below so the formula text matches.
$ws.Range("E6").Formula = "=E3*3+E9*2+10"

# Row 2 ("Motor: NPC 2212"): cost formula zeroed out (motor no longer costed
# in the running total the same way).
$ws.Range("F2").Formula = "=0"

# Remove the "panel retainer 5jb82" row entirely - rows below it (Panel
# large, panel small, Mcmaster wheels, ...) shift up by one.
$ws.Rows("10:10").Delete()

# Match the author's final cell selection.
$null = $ws.Range("H8").Select()
